# Report Screenshot dan Video
# Append 65 more saved-tag rows (rows 26-90) to column A of the "TagSave"
# sheet (the active sheet when this runs is "Login", so the sheet is
# looked up explicitly by name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TagSave")

$values = @(
    "YsIgF2cHnN", "RpPuTMRrFJ", "PxYGyOXlue", "T9Wv2Z21HX", "xzvugZSjnu", "O2XQE27noP", "vBTO68XlLc", "92NdKdDuI3",
    "KXWNFDQcAU", "Lv74x8QNgF", "EgH1BXMdDt", "syiS76rjYp", "82b6G6ZTWp", "HVSNS4kXHA", "qz6C02DmUl", "787ZL9LJDV",
    "1aPf1qDOtn", "O9EZI3m8qq", "vFnRaZAJtp", "IfCgdGe1H7", "w4qQieb4tQ", "SVfas0iNyK", "ri9yfS50Du", "R8B4EZHsnV",
    "OrZon85kC6", "4aL5r6uq1N", "k07l5vTt6g", "QP1bTr7fYB", "0zYI61h88z", "8iV4KPBO2P", "dnUUl9ywom", "qGNqugSxOt",
    "Fy5iKPrAuY", "0uSpYmRgO6", "dJCIYuyY0r", "YXUQn7gJcB", "IpzsfX8eOq", "6pNRFcNwxu", "Aaw9I1L4ig", "wyMLuBf7iO",
    "TDAirExIr1", "NrhPXZJ1Fi", "JaAWp73dve", "3O5VxPDjtl", "y6ZJJm764k", "YbLavVR7mV", "PjZ8C0jHym", "j6JA2UK2cL",
    "V2uirQgFW2", "aVT79aSw7b", "7wPuUBbQVd", "QBhqOPAgmc", "K2H9nFYon4", "ilk76xav59", "KE0v7gCFJZ", "ldfTAOYBFq",
    "wFMxzgPukX", "b9hvcO02eP", "KpNXOENfsB", "1mMSths9W4", "7aEDQTJb58", "kS3Bn2kx6u", "2fKKWhfFCv", "mOmT4UJzHX",
    "97MUUlFJJZ"
)

$startRow = 26
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

